# Replace the hard-coded "Jakarta" sign-off location with the
# ${authority_location} merge placeholder in every "Jakarta, ${docs_date}"
# style paragraph throughout the document, leaving the trailing
# comma/space(s) and the run's formatting untouched.
$d = $word.ActiveDocument

$rng = $d.Content
$rng.Find.ClearFormatting()
$rng.Find.Replacement.ClearFormatting()
[void]$rng.Find.Execute("Jakarta", $true, $false, $false, $false, $false, $true, 1, $false, '${authority_location}', 2)
